$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows for 2025-02-08 ---
$ws.Range("D23").Value = $true
$ws.Range("C24").Value = $true
$ws.Range("D24").Value = $true

# --- Append new rows for 2025-02-09 ---
# Force column A to text first so the "YYYY-MM-DD" strings are stored as
# plain text (matching the rest of the sheet) instead of being
# auto-converted into date serial numbers, then restore the default
# (unstyled) look once the values are in place.
$dateRange = $ws.Range("A26:A28")
$dateRange.NumberFormat = "@"

$ws.Cells.Item(26, 1).Value = "2025-02-09"
$ws.Cells.Item(26, 2).Value = "sleep"
$ws.Cells.Item(26, 3).Value = $false
$ws.Cells.Item(26, 4).Value = $false

$ws.Cells.Item(27, 1).Value = "2025-02-09"
$ws.Cells.Item(27, 2).Value = "activity"
$ws.Cells.Item(27, 3).Value = $true
$ws.Cells.Item(27, 4).Value = $true

$ws.Cells.Item(28, 1).Value = "2025-02-09"
$ws.Cells.Item(28, 2).Value = "weekly_activity"
$ws.Cells.Item(28, 3).Value = $false
$ws.Cells.Item(28, 4).Value = $false

$dateRange.Style = "Normal"
